$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would be misread as a number (trailing
# zeros like "1.00" / "37.20" collapse to "1" / "37.2" under General
# auto-detection) -- force Text format first so they stay literal strings.
$textCells = "D4", "D14", "D28", "D32", "D42", "D47", "D50", "D51"
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '68.420.25'
$ws.Range("E2").Value = '  -1.52%  '

$ws.Range("D3").Value = '3.839.55'
$ws.Range("E3").Value = '  -1.10%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '601.24'
$ws.Range("E5").Value = '  -0.40%  '

$ws.Range("D6").Value = '170.03'
$ws.Range("E6").Value = '  +0.39%  '

$ws.Range("D7").Value = '3.838.96'
$ws.Range("E7").Value = '  -1.17%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").Value = '0.527'
$ws.Range("E9").Value = '  -1.26%  '

$ws.Range("E10").Value = '  -1.40%  '

$ws.Range("E11").Value = '  +1.91%  '

$ws.Range("E12").Value = '  -1.98%  '

$ws.Range("E13").Value = '  +3.74%  '

$ws.Range("D14").Value = '37.20'
$ws.Range("E14").Value = '  -2.61%  '

$ws.Range("D15").Value = '4.482.19'
$ws.Range("E15").Value = '  -1.21%  '

$ws.Range("D16").Value = '3.836.14'
$ws.Range("E16").Value = '  -1.26%  '

$ws.Range("D17").Value = '68.400.86'
$ws.Range("E17").Value = '  -1.57%  '

$ws.Range("D18").Value = '18.55'
$ws.Range("E18").Value = '  -1.25%  '

$ws.Range("E19").Value = '  -2.83%  '

$ws.Range("E20").Value = '  -0.65%  '

$ws.Range("D21").Value = '11.12'
$ws.Range("E21").Value = '  -0.39%  '

$ws.Range("D22").Value = '470.08'
$ws.Range("E22").Value = '  -3.97%  '

$ws.Range("D23").Value = '0.738'
$ws.Range("E23").Value = '  -1.20%  '

$ws.Range("E24").Value = '  -3.03%  '

$ws.Range("D25").Value = '83.26'
$ws.Range("E25").Value = '  -2.15%  '

$ws.Range("E26").Value = '  -2.20%  '

$ws.Range("E27").Value = '  -1.12%  '

$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.07%  '

$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '10.01'
$ws.Range("E29").Value = '  -1.26%  '

$ws.Range("E30").Value = '  -0.17%  '

$ws.Range("D31").Value = '3.988.06'
$ws.Range("E31").Value = '  -1.13%  '

$ws.Range("D32").Value = '7.70'
$ws.Range("E32").Value = '  -1.06%  '

$ws.Range("E33").Value = '  -0.75%  '

$ws.Range("E34").Value = '  -4.30%  '

$ws.Range("E35").Value = '  -0.65%  '

$ws.Range("D36").Value = '3.802.80'
$ws.Range("E36").Value = '  -1.15%  '

$ws.Range("E37").Value = '  -1.96%  '

$ws.Range("D38").Value = '3.68'
$ws.Range("E38").Value = '  +12.31%  '

$ws.Range("E39").Value = '  -1.18%  '

$ws.Range("E40").Value = '  -1.02%  '

$ws.Range("D41").Value = '5.96'
$ws.Range("E41").Value = '  -2.38%  '

$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.06%  '

$ws.Range("E43").Value = '  -3.12%  '

$ws.Range("D44").Value = '1.99'
$ws.Range("E44").Value = '  -5.49%  '

$ws.Range("D45").Value = '8.82'
$ws.Range("E45").Value = '  +1.50%  '

$ws.Range("D46").Value = '419.26'
$ws.Range("E46").Value = '  -3.98%  '

$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.00%  '

$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").Value = '0.000293'
$ws.Range("E48").Value = '  +8.14%  '

$ws.Range("E49").Value = '  -2.06%  '

$ws.Range("D50").Value = '26.20'
$ws.Range("E50").Value = '  +4.78%  '

$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = '0.0360'
$ws.Range("E51").Value = '  -2.50%  '
